# Reporte Mids - remove columns that are no longer reported
# (ACTIVO, USUARIO, FECHA_ACTUAL) now that an exclusion list is
# maintained separately.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Datos")

# Delete whole columns from right to left so the earlier column
# letters (G, S) still point at the right columns when deleted.
$ws.Columns("Y").Delete()   # FECHA_ACTUAL
$ws.Columns("S").Delete()   # USUARIO
$ws.Columns("G").Delete()   # ACTIVO
